$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the columns that get rearranged (D, L, M, N, O, P, Q, R, S, T)
# for every data row (2..17) before writing anything, so reads never see
# a partially-updated sheet.

$D2  = $ws.Range("D2").Value2
$L2  = $ws.Range("L2").Value2
$M2  = $ws.Range("M2").Value2
$N2  = $ws.Range("N2").Value2
$O2  = $ws.Range("O2").Value2
$P2  = $ws.Range("P2").Value2
$Q2  = $ws.Range("Q2").Value2
$R2  = $ws.Range("R2").Value2
$S2  = $ws.Range("S2").Value2
$T2  = $ws.Range("T2").Value2

$D3  = $ws.Range("D3").Value2
$L3  = $ws.Range("L3").Value2
$M3  = $ws.Range("M3").Value2
$N3  = $ws.Range("N3").Value2
$O3  = $ws.Range("O3").Value2
$P3  = $ws.Range("P3").Value2
$Q3  = $ws.Range("Q3").Value2
$R3  = $ws.Range("R3").Value2
$S3  = $ws.Range("S3").Value2
$T3  = $ws.Range("T3").Value2

$D4  = $ws.Range("D4").Value2
$L4  = $ws.Range("L4").Value2
$M4  = $ws.Range("M4").Value2
$N4  = $ws.Range("N4").Value2
$O4  = $ws.Range("O4").Value2
$P4  = $ws.Range("P4").Value2
$Q4  = $ws.Range("Q4").Value2
$R4  = $ws.Range("R4").Value2
$S4  = $ws.Range("S4").Value2
$T4  = $ws.Range("T4").Value2

$D5  = $ws.Range("D5").Value2
$L5  = $ws.Range("L5").Value2
$M5  = $ws.Range("M5").Value2
$N5  = $ws.Range("N5").Value2
$O5  = $ws.Range("O5").Value2
$P5  = $ws.Range("P5").Value2
$Q5  = $ws.Range("Q5").Value2
$R5  = $ws.Range("R5").Value2
$S5  = $ws.Range("S5").Value2
$T5  = $ws.Range("T5").Value2

$D6  = $ws.Range("D6").Value2
$L6  = $ws.Range("L6").Value2
$M6  = $ws.Range("M6").Value2
$N6  = $ws.Range("N6").Value2
$O6  = $ws.Range("O6").Value2
$P6  = $ws.Range("P6").Value2
$Q6  = $ws.Range("Q6").Value2
$R6  = $ws.Range("R6").Value2
$S6  = $ws.Range("S6").Value2
$T6  = $ws.Range("T6").Value2

$D7  = $ws.Range("D7").Value2
$L7  = $ws.Range("L7").Value2
$M7  = $ws.Range("M7").Value2
$N7  = $ws.Range("N7").Value2
$O7  = $ws.Range("O7").Value2
$P7  = $ws.Range("P7").Value2
$Q7  = $ws.Range("Q7").Value2
$R7  = $ws.Range("R7").Value2
$S7  = $ws.Range("S7").Value2
$T7  = $ws.Range("T7").Value2

$D8  = $ws.Range("D8").Value2
$L8  = $ws.Range("L8").Value2
$M8  = $ws.Range("M8").Value2
$N8  = $ws.Range("N8").Value2
$O8  = $ws.Range("O8").Value2
$P8  = $ws.Range("P8").Value2
$Q8  = $ws.Range("Q8").Value2
$R8  = $ws.Range("R8").Value2
$S8  = $ws.Range("S8").Value2
$T8  = $ws.Range("T8").Value2

$D9  = $ws.Range("D9").Value2
$L9  = $ws.Range("L9").Value2
$M9  = $ws.Range("M9").Value2
$N9  = $ws.Range("N9").Value2
$O9  = $ws.Range("O9").Value2
$P9  = $ws.Range("P9").Value2
$Q9  = $ws.Range("Q9").Value2
$R9  = $ws.Range("R9").Value2
$S9  = $ws.Range("S9").Value2
$T9  = $ws.Range("T9").Value2

$D10 = $ws.Range("D10").Value2
$L10 = $ws.Range("L10").Value2
$M10 = $ws.Range("M10").Value2
$N10 = $ws.Range("N10").Value2
$O10 = $ws.Range("O10").Value2
$P10 = $ws.Range("P10").Value2
$Q10 = $ws.Range("Q10").Value2
$R10 = $ws.Range("R10").Value2
$S10 = $ws.Range("S10").Value2
$T10 = $ws.Range("T10").Value2

$D11 = $ws.Range("D11").Value2
$L11 = $ws.Range("L11").Value2
$M11 = $ws.Range("M11").Value2
$N11 = $ws.Range("N11").Value2
$O11 = $ws.Range("O11").Value2
$P11 = $ws.Range("P11").Value2
$Q11 = $ws.Range("Q11").Value2
$R11 = $ws.Range("R11").Value2
$S11 = $ws.Range("S11").Value2
$T11 = $ws.Range("T11").Value2

$D12 = $ws.Range("D12").Value2
$L12 = $ws.Range("L12").Value2
$M12 = $ws.Range("M12").Value2
$N12 = $ws.Range("N12").Value2
$O12 = $ws.Range("O12").Value2
$P12 = $ws.Range("P12").Value2
$Q12 = $ws.Range("Q12").Value2
$R12 = $ws.Range("R12").Value2
$S12 = $ws.Range("S12").Value2
$T12 = $ws.Range("T12").Value2

$D13 = $ws.Range("D13").Value2
$L13 = $ws.Range("L13").Value2
$M13 = $ws.Range("M13").Value2
$N13 = $ws.Range("N13").Value2
$O13 = $ws.Range("O13").Value2
$P13 = $ws.Range("P13").Value2
$Q13 = $ws.Range("Q13").Value2
$R13 = $ws.Range("R13").Value2
$S13 = $ws.Range("S13").Value2
$T13 = $ws.Range("T13").Value2

$D14 = $ws.Range("D14").Value2
$L14 = $ws.Range("L14").Value2
$M14 = $ws.Range("M14").Value2
$N14 = $ws.Range("N14").Value2
$O14 = $ws.Range("O14").Value2
$P14 = $ws.Range("P14").Value2
$Q14 = $ws.Range("Q14").Value2
$R14 = $ws.Range("R14").Value2
$S14 = $ws.Range("S14").Value2
$T14 = $ws.Range("T14").Value2

$D15 = $ws.Range("D15").Value2
$L15 = $ws.Range("L15").Value2
$M15 = $ws.Range("M15").Value2
$N15 = $ws.Range("N15").Value2
$O15 = $ws.Range("O15").Value2
$P15 = $ws.Range("P15").Value2
$Q15 = $ws.Range("Q15").Value2
$R15 = $ws.Range("R15").Value2
$S15 = $ws.Range("S15").Value2
$T15 = $ws.Range("T15").Value2

$D16 = $ws.Range("D16").Value2
$L16 = $ws.Range("L16").Value2
$M16 = $ws.Range("M16").Value2
$N16 = $ws.Range("N16").Value2
$O16 = $ws.Range("O16").Value2
$P16 = $ws.Range("P16").Value2
$Q16 = $ws.Range("Q16").Value2
$R16 = $ws.Range("R16").Value2
$S16 = $ws.Range("S16").Value2
$T16 = $ws.Range("T16").Value2

$D17 = $ws.Range("D17").Value2
$L17 = $ws.Range("L17").Value2
$M17 = $ws.Range("M17").Value2
$N17 = $ws.Range("N17").Value2
$O17 = $ws.Range("O17").Value2
$P17 = $ws.Range("P17").Value2
$Q17 = $ws.Range("Q17").Value2
$R17 = $ws.Range("R17").Value2
$S17 = $ws.Range("S17").Value2
$T17 = $ws.Range("T17").Value2

# Write row 2 <- source row 8
$ws.Range("D2").Value  = $D8
$ws.Range("L2").Value  = $L8
$ws.Range("M2").Value  = $M8
$ws.Range("N2").Value  = $N8
$ws.Range("O2").Value  = $O8
$ws.Range("P2").Value  = $P8
$ws.Range("Q2").Value  = $Q8
$ws.Range("R2").Value  = $R8
$ws.Range("S2").Value  = $S8
$ws.Range("T2").Value  = $T8

# Write row 3 <- source row 10
$ws.Range("D3").Value  = $D10
$ws.Range("L3").Value  = $L10
$ws.Range("M3").Value  = $M10
$ws.Range("N3").Value  = $N10
$ws.Range("O3").Value  = $O10
$ws.Range("P3").Value  = $P10
$ws.Range("Q3").Value  = $Q10
$ws.Range("R3").Value  = $R10
$ws.Range("S3").Value  = $S10
$ws.Range("T3").Value  = $T10

# Write row 4 <- source row 12
$ws.Range("D4").Value  = $D12
$ws.Range("L4").Value  = $L12
$ws.Range("M4").Value  = $M12
$ws.Range("N4").Value  = $N12
$ws.Range("O4").Value  = $O12
$ws.Range("P4").Value  = $P12
$ws.Range("Q4").Value  = $Q12
$ws.Range("R4").Value  = $R12
$ws.Range("S4").Value  = $S12
$ws.Range("T4").Value  = $T12

# Write row 5 <- source row 14
$ws.Range("D5").Value  = $D14
$ws.Range("L5").Value  = $L14
$ws.Range("M5").Value  = $M14
$ws.Range("N5").Value  = $N14
$ws.Range("O5").Value  = $O14
$ws.Range("P5").Value  = $P14
$ws.Range("Q5").Value  = $Q14
$ws.Range("R5").Value  = $R14
$ws.Range("S5").Value  = $S14
$ws.Range("T5").Value  = $T14

# Write row 6 <- source row 11
$ws.Range("D6").Value  = $D11
$ws.Range("L6").Value  = $L11
$ws.Range("M6").Value  = $M11
$ws.Range("N6").Value  = $N11
$ws.Range("O6").Value  = $O11
$ws.Range("P6").Value  = $P11
$ws.Range("Q6").Value  = $Q11
$ws.Range("R6").Value  = $R11
$ws.Range("S6").Value  = $S11
$ws.Range("T6").Value  = $T11

# Write row 7 <- source row 9
$ws.Range("D7").Value  = $D9
$ws.Range("L7").Value  = $L9
$ws.Range("M7").Value  = $M9
$ws.Range("N7").Value  = $N9
$ws.Range("O7").Value  = $O9
$ws.Range("P7").Value  = $P9
$ws.Range("Q7").Value  = $Q9
$ws.Range("R7").Value  = $R9
$ws.Range("S7").Value  = $S9
$ws.Range("T7").Value  = $T9

# Write row 8 <- source row 17
$ws.Range("D8").Value  = $D17
$ws.Range("L8").Value  = $L17
$ws.Range("M8").Value  = $M17
$ws.Range("N8").Value  = $N17
$ws.Range("O8").Value  = $O17
$ws.Range("P8").Value  = $P17
$ws.Range("Q8").Value  = $Q17
$ws.Range("R8").Value  = $R17
$ws.Range("S8").Value  = $S17
$ws.Range("T8").Value  = $T17

# Write row 9 <- source row 13
$ws.Range("D9").Value  = $D13
$ws.Range("L9").Value  = $L13
$ws.Range("M9").Value  = $M13
$ws.Range("N9").Value  = $N13
$ws.Range("O9").Value  = $O13
$ws.Range("P9").Value  = $P13
$ws.Range("Q9").Value  = $Q13
$ws.Range("R9").Value  = $R13
$ws.Range("S9").Value  = $S13
$ws.Range("T9").Value  = $T13

# Write row 10 <- source row 6
$ws.Range("D10").Value = $D6
$ws.Range("L10").Value = $L6
$ws.Range("M10").Value = $M6
$ws.Range("N10").Value = $N6
$ws.Range("O10").Value = $O6
$ws.Range("P10").Value = $P6
$ws.Range("Q10").Value = $Q6
$ws.Range("R10").Value = $R6
$ws.Range("S10").Value = $S6
$ws.Range("T10").Value = $T6

# Write row 11 <- source row 15
$ws.Range("D11").Value = $D15
$ws.Range("L11").Value = $L15
$ws.Range("M11").Value = $M15
$ws.Range("N11").Value = $N15
$ws.Range("O11").Value = $O15
$ws.Range("P11").Value = $P15
$ws.Range("Q11").Value = $Q15
$ws.Range("R11").Value = $R15
$ws.Range("S11").Value = $S15
$ws.Range("T11").Value = $T15

# Write row 12 <- source row 16
$ws.Range("D12").Value = $D16
$ws.Range("L12").Value = $L16
$ws.Range("M12").Value = $M16
$ws.Range("N12").Value = $N16
$ws.Range("O12").Value = $O16
$ws.Range("P12").Value = $P16
$ws.Range("Q12").Value = $Q16
$ws.Range("R12").Value = $R16
$ws.Range("S12").Value = $S16
$ws.Range("T12").Value = $T16

# Write row 13 <- source row 5
$ws.Range("D13").Value = $D5
$ws.Range("L13").Value = $L5
$ws.Range("M13").Value = $M5
$ws.Range("N13").Value = $N5
$ws.Range("O13").Value = $O5
$ws.Range("P13").Value = $P5
$ws.Range("Q13").Value = $Q5
$ws.Range("R13").Value = $R5
$ws.Range("S13").Value = $S5
$ws.Range("T13").Value = $T5

# Write row 14 <- source row 2
$ws.Range("D14").Value = $D2
$ws.Range("L14").Value = $L2
$ws.Range("M14").Value = $M2
$ws.Range("N14").Value = $N2
$ws.Range("O14").Value = $O2
$ws.Range("P14").Value = $P2
$ws.Range("Q14").Value = $Q2
$ws.Range("R14").Value = $R2
$ws.Range("S14").Value = $S2
$ws.Range("T14").Value = $T2

# Write row 15 <- source row 3
$ws.Range("D15").Value = $D3
$ws.Range("L15").Value = $L3
$ws.Range("M15").Value = $M3
$ws.Range("N15").Value = $N3
$ws.Range("O15").Value = $O3
$ws.Range("P15").Value = $P3
$ws.Range("Q15").Value = $Q3
$ws.Range("R15").Value = $R3
$ws.Range("S15").Value = $S3
$ws.Range("T15").Value = $T3

# Write row 16 <- source row 4
$ws.Range("D16").Value = $D4
$ws.Range("L16").Value = $L4
$ws.Range("M16").Value = $M4
$ws.Range("N16").Value = $N4
$ws.Range("O16").Value = $O4
$ws.Range("P16").Value = $P4
$ws.Range("Q16").Value = $Q4
$ws.Range("R16").Value = $R4
$ws.Range("S16").Value = $S4
$ws.Range("T16").Value = $T4

# Write row 17 <- source row 7
$ws.Range("D17").Value = $D7
$ws.Range("L17").Value = $L7
$ws.Range("M17").Value = $M7
$ws.Range("N17").Value = $N7
$ws.Range("O17").Value = $O7
$ws.Range("P17").Value = $P7
$ws.Range("Q17").Value = $Q7
$ws.Range("R17").Value = $R7
$ws.Range("S17").Value = $S7
$ws.Range("T17").Value = $T7
